$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force columns A and E to text format so large integer IDs are preserved exactly
$ws.Range("A2:A18").NumberFormat = "@"
$ws.Range("E2:E18").NumberFormat = "@"

$ws.Cells.Item(2, 1).Value = "4529285391522200320"
$ws.Cells.Item(2, 2).Value = 275.9254491025025
$ws.Cells.Item(2, 3).Value = 21.76870525491716
$ws.Cells.Item(2, 4).Value = 3.389339208602905
$ws.Cells.Item(2, 5).Value = "4529285391531266304"
$ws.Cells.Item(2, 6).Value = 3.39525294303894
$ws.Cells.Item(2, 7).Value = 0.0001677126147536809

$ws.Cells.Item(3, 1).Value = "5426587107145955712"
$ws.Cells.Item(3, 2).Value = 142.6736805260494
$ws.Cells.Item(3, 3).Value = -40.46642826305888
$ws.Cells.Item(3, 4).Value = 3.674772262573242
$ws.Cells.Item(3, 5).Value = "5426587107149861120"
$ws.Cells.Item(3, 6).Value = 4.803853511810303
$ws.Cells.Item(3, 7).Value = 0.0002820617191807018

$ws.Cells.Item(4, 1).Value = "4432032213656433536"
$ws.Cells.Item(4, 2).Value = 247.7282713167697
$ws.Cells.Item(4, 3).Value = 1.983620889639908
$ws.Cells.Item(4, 4).Value = 4.066857814788818
$ws.Cells.Item(4, 5).Value = "4432032213656434176"
$ws.Cells.Item(4, 6).Value = 6.452131748199463
$ws.Cells.Item(4, 7).Value = 0.0003843007249245102

$ws.Cells.Item(5, 1).Value = "6541802406664428672"
$ws.Cells.Item(5, 2).Value = 346.7194246886781
$ws.Cells.Item(5, 3).Value = -43.52041179748344
$ws.Cells.Item(5, 4).Value = 4.234705448150635
$ws.Cells.Item(5, 5).Value = "6541802402371581568"
$ws.Cells.Item(5, 6).Value = 6.503548145294189
$ws.Cells.Item(5, 7).Value = 0.0004380720153890083

$ws.Cells.Item(6, 1).Value = "6719152945029845376"
$ws.Cells.Item(6, 2).Value = 286.6054157588619
$ws.Cells.Item(6, 3).Value = -37.06488085942344
$ws.Cells.Item(6, 4).Value = 4.717123985290527
$ws.Cells.Item(6, 5).Value = "6719152945032456832"
$ws.Cells.Item(6, 6).Value = 4.760257244110107
$ws.Cells.Item(6, 7).Value = 0.0003869906325269605

$ws.Cells.Item(7, 1).Value = "6719152945032456832"
$ws.Cells.Item(7, 2).Value = 286.6053079615535
$ws.Cells.Item(7, 3).Value = -37.06450354954132
$ws.Cells.Item(7, 4).Value = 4.760257244110107
$ws.Cells.Item(7, 5).Value = "6719152945029845376"
$ws.Cells.Item(7, 6).Value = 4.717123985290527
$ws.Cells.Item(7, 7).Value = 0.0003869906325269605

$ws.Cells.Item(8, 1).Value = "4343066192373234048"
$ws.Cells.Item(8, 2).Value = 241.0920131228604
$ws.Cells.Item(8, 3).Value = -11.37304460452119
$ws.Cells.Item(8, 4).Value = 4.767579555511475
$ws.Cells.Item(8, 5).Value = "4343066192367555200"
$ws.Cells.Item(8, 6).Value = 4.771558284759521
$ws.Cells.Item(8, 7).Value = 0.0003030379884824618

$ws.Cells.Item(9, 1).Value = "4343066192367555200"
$ws.Cells.Item(9, 2).Value = 241.0919870016956
$ws.Cells.Item(9, 3).Value = -11.37334655855848
$ws.Cells.Item(9, 4).Value = 4.771558284759521
$ws.Cells.Item(9, 5).Value = "4343066192373234048"
$ws.Cells.Item(9, 6).Value = 4.767579555511475
$ws.Cells.Item(9, 7).Value = 0.0003030379884824618

$ws.Cells.Item(10, 1).Value = "6860945174279114880"
$ws.Cells.Item(10, 2).Value = 307.2149851888128
$ws.Cells.Item(10, 3).Value = -17.81372857084146
$ws.Cells.Item(10, 4).Value = 4.794205188751221
$ws.Cells.Item(10, 5).Value = "6860945174275852416"
$ws.Cells.Item(10, 6).Value = 6.696762561798096
$ws.Cells.Item(10, 7).Value = 0.0004518621805492849

$ws.Cells.Item(11, 1).Value = "5426587107149861120"
$ws.Cells.Item(11, 2).Value = 142.6740032632142
$ws.Cells.Item(11, 3).Value = -40.466567086598
$ws.Cells.Item(11, 4).Value = 4.803853511810303
$ws.Cells.Item(11, 5).Value = "5426587107145955712"
$ws.Cells.Item(11, 6).Value = 3.674772262573242
$ws.Cells.Item(11, 7).Value = 0.0002820617191807018

$ws.Cells.Item(12, 1).Value = "657244586015485440"
$ws.Cells.Item(12, 2).Value = 123.0535363160312
$ws.Cells.Item(12, 3).Value = 17.64700231993405
$ws.Cells.Item(12, 4).Value = 5.410910606384277
$ws.Cells.Item(12, 5).Value = "657244521593509376"
$ws.Cells.Item(12, 6).Value = 5.792698383331299
$ws.Cells.Item(12, 7).Value = 0.0003145394475308107

$ws.Cells.Item(13, 1).Value = "6724105656508792576"
$ws.Cells.Item(13, 2).Value = 271.7079518856002
$ws.Cells.Item(13, 3).Value = -43.42567811693711
$ws.Cells.Item(13, 4).Value = 5.581056118011475
$ws.Cells.Item(13, 5).Value = "6724105660828668032"
$ws.Cells.Item(13, 6).Value = 5.644303798675537
$ws.Cells.Item(13, 7).Value = 0.0004930990478208497

$ws.Cells.Item(14, 1).Value = "6724105660828668032"
$ws.Cells.Item(14, 2).Value = 271.7079719599973
$ws.Cells.Item(14, 3).Value = -43.42518523347083
$ws.Cells.Item(14, 4).Value = 5.644303798675537
$ws.Cells.Item(14, 5).Value = "6724105656508792576"
$ws.Cells.Item(14, 6).Value = 5.581056118011475
$ws.Cells.Item(14, 7).Value = 0.0004930990478208497

$ws.Cells.Item(15, 1).Value = "2803936143261926272"
$ws.Cells.Item(15, 2).Value = 13.74271321842037
$ws.Cells.Item(15, 3).Value = 23.62814488153855
$ws.Cells.Item(15, 4).Value = 5.725862503051758
$ws.Cells.Item(15, 5).Value = "2803936138967498624"
$ws.Cells.Item(15, 6).Value = 6.123550891876221
$ws.Cells.Item(15, 7).Value = 0.0003119170200913206

$ws.Cells.Item(16, 1).Value = "4373199678620639744"
$ws.Cells.Item(16, 2).Value = 262.5984522008857
$ws.Cells.Item(16, 3).Value = -1.063482522283796
$ws.Cells.Item(16, 4).Value = 5.779087543487549
$ws.Cells.Item(16, 5).Value = "4373199682919087616"
$ws.Cells.Item(16, 6).Value = 5.801146030426025
$ws.Cells.Item(16, 7).Value = 0.000199786113763143

$ws.Cells.Item(17, 1).Value = "657244521593509376"
$ws.Cells.Item(17, 2).Value = 123.0536499059178
$ws.Cells.Item(17, 3).Value = 17.64729764719009
$ws.Cells.Item(17, 4).Value = 5.792698383331299
$ws.Cells.Item(17, 5).Value = "657244586015485440"
$ws.Cells.Item(17, 6).Value = 5.410910606384277
$ws.Cells.Item(17, 7).Value = 0.0003145394475308107

$ws.Cells.Item(18, 1).Value = "4373199682919087616"
$ws.Cells.Item(18, 2).Value = 262.5985706282946
$ws.Cells.Item(18, 3).Value = -1.063643439256715
$ws.Cells.Item(18, 4).Value = 5.801146030426025
$ws.Cells.Item(18, 5).Value = "4373199678620639744"
$ws.Cells.Item(18, 6).Value = 5.779087543487549
$ws.Cells.Item(18, 7).Value = 0.000199786113763143

# Remove the temporary text-number-format styling so cells keep the default style,
# while preserving the text values already written.
$ws.Range("A2:A18").ClearFormats()
$ws.Range("E2:E18").ClearFormats()
